$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 163: insert AR163 ---
$ws.Range("AR163").Value = 15.4762831

# --- Row 192: corrected values ---
$ws.Range("G192").Value = 19.4216679
$ws.Range("AI192").Value = 9.909747299999999
$ws.Range("AW192").Value = 29.0146201
$ws.Range("BB192").Value = 17.0977809

# --- Row 193 ---
$ws.Range("B193").Value = 19.6052632
$ws.Range("C193").Value = 34.2661769
$ws.Range("D193").Value = 29.2593466
$ws.Range("F193").Value = 24.9705315
$ws.Range("G193").Value = 19.1479942
$ws.Range("H193").Value = 16.7845788
$ws.Range("I193").Value = 9.848635700000001
$ws.Range("J193").Value = 13.6963696
$ws.Range("K193").Value = 15.5461824
$ws.Range("L193").Value = 25.927384
$ws.Range("M193").Value = 31.3099781
$ws.Range("O193").Value = 15.1958081
$ws.Range("P193").Value = 27.3616873
$ws.Range("Q193").Value = 30.0733337
$ws.Range("R193").Value = 20.1294749
$ws.Range("S193").Value = 24.8541971
$ws.Range("T193").Value = 25.5222509
$ws.Range("U193").Value = 26.3444579
$ws.Range("V193").Value = 33.3808065
$ws.Range("W193").Value = 10.9760642
$ws.Range("X193").Value = 14.8123237
$ws.Range("Y193").Value = 10.47994
$ws.Range("Z193").Value = 16.5274594
$ws.Range("AA193").Value = 20.1939928
$ws.Range("AB193").Value = 25.0355124
$ws.Range("AD193").Value = 36.1420749
$ws.Range("AE193").Value = 25.1103446
$ws.Range("AF193").Value = 20.5472858
$ws.Range("AG193").Value = 24.6669186
$ws.Range("AH193").Value = 26.8841593
$ws.Range("AI193").Value = 10.5946556
$ws.Range("AJ193").Value = 11.7983444
$ws.Range("AK193").Value = 18.216194
$ws.Range("AL193").Value = 25.2037305
$ws.Range("AM193").Value = 11.0392656
$ws.Range("AN193").Value = 21.4413704
$ws.Range("AO193").Value = 28.8569511
$ws.Range("AP193").Value = 14.9500496
$ws.Range("AQ193").Value = 14.6292209
$ws.Range("AS193").Value = 12.956355
$ws.Range("AT193").Value = 28.2340734
$ws.Range("AU193").Value = 22.967095
$ws.Range("AV193").Value = 29.4090157
$ws.Range("AW193").Value = 28.8600958
$ws.Range("AX193").Value = 25.8212266
$ws.Range("AY193").Value = 18.2934525
$ws.Range("BA193").Value = 8.7250996
$ws.Range("BB193").Value = 17.4668029
$ws.Range("BC193").Value = 20.5820599
$ws.Range("BD193").Value = 20.430414
$ws.Range("BE193").Value = 18.9908921

# --- Row 194 ---
$ws.Range("A194").Value = "11 08 2020"
$ws.Range("B194").Value = 19.6826758
$ws.Range("C194").Value = 33.8765221
$ws.Range("D194").Value = 28.32004
$ws.Range("F194").Value = 24.23096
$ws.Range("G194").Value = 18.9054907
$ws.Range("H194").Value = 17.0580532
$ws.Range("I194").Value = 9.935961600000001
$ws.Range("J194").Value = 13.3215962
$ws.Range("K194").Value = 14.571263
$ws.Range("L194").Value = 25.7175926
$ws.Range("M194").Value = 31.0010299
$ws.Range("O194").Value = 13.9184397
$ws.Range("P194").Value = 27.1160785
$ws.Range("Q194").Value = 30.3600343
$ws.Range("R194").Value = 19.982742
$ws.Range("S194").Value = 25.2478949
$ws.Range("T194").Value = 24.6751869
$ws.Range("U194").Value = 25.528797
$ws.Range("V194").Value = 32.7764527
$ws.Range("W194").Value = 11.0868691
$ws.Range("X194").Value = 14.3528896
$ws.Range("Y194").Value = 10.5667078
$ws.Range("Z194").Value = 16.7785443
$ws.Range("AA194").Value = 20.3893285
$ws.Range("AB194").Value = 24.7302256
$ws.Range("AD194").Value = 35.2059278
$ws.Range("AE194").Value = 24.0510772
$ws.Range("AF194").Value = 20.7663675
$ws.Range("AG194").Value = 24.7432755
$ws.Range("AH194").Value = 26.5367746
$ws.Range("AI194").Value = 10.1242443
$ws.Range("AJ194").Value = 11.9191606
$ws.Range("AK194").Value = 18.5240634
$ws.Range("AL194").Value = 25.3750702
$ws.Range("AM194").Value = 11.1327329
$ws.Range("AN194").Value = 21.0373579
$ws.Range("AO194").Value = 28.9029094
$ws.Range("AP194").Value = 14.8865666
$ws.Range("AQ194").Value = 14.5729544
$ws.Range("AS194").Value = 11.8461968
$ws.Range("AT194").Value = 27.6780118
$ws.Range("AU194").Value = 24.8410146
$ws.Range("AV194").Value = 29.4496858
$ws.Range("AW194").Value = 28.3820293
$ws.Range("AX194").Value = 25.6661444
$ws.Range("AY194").Value = 18.3346009
$ws.Range("BA194").Value = 10.1078282
$ws.Range("BB194").Value = 17.3257563
$ws.Range("BC194").Value = 20.3933499
$ws.Range("BD194").Value = 21.64249
$ws.Range("BE194").Value = 19.2265712

# --- Row 195 ---
$ws.Range("A195").Value = "12 08 2020"
$ws.Range("B195").Value = 19.4235589
$ws.Range("C195").Value = 33.9847951
$ws.Range("D195").Value = 28.1400731
$ws.Range("F195").Value = 23.9119698
$ws.Range("G195").Value = 18.7414373
$ws.Range("H195").Value = 16.7848279
$ws.Range("I195").Value = 10.3551438
$ws.Range("J195").Value = 13.4439359
$ws.Range("K195").Value = 13.8987609
$ws.Range("L195").Value = 25.2462839
$ws.Range("M195").Value = 30.46785
$ws.Range("O195").Value = 15
$ws.Range("P195").Value = 26.4654174
$ws.Range("Q195").Value = 30.4311724
$ws.Range("R195").Value = 19.3729531
$ws.Range("S195").Value = 24.7837624
$ws.Range("T195").Value = 24.6845356
$ws.Range("U195").Value = 26.6649871
$ws.Range("V195").Value = 33.0114805
$ws.Range("W195").Value = 10.9608671
$ws.Range("X195").Value = 14.1634027
$ws.Range("Y195").Value = 11.0333081
$ws.Range("Z195").Value = 17.0224811
$ws.Range("AA195").Value = 20.0185566
$ws.Range("AB195").Value = 25.3678791
$ws.Range("AD195").Value = 35.7270295
$ws.Range("AE195").Value = 24.1301568
$ws.Range("AF195").Value = 20.7941167
$ws.Range("AG195").Value = 26.345527
$ws.Range("AH195").Value = 25.7799169
$ws.Range("AI195").Value = 10.3171496
$ws.Range("AJ195").Value = 11.8360585
$ws.Range("AK195").Value = 18.0696453
$ws.Range("AL195").Value = 24.2451588
$ws.Range("AM195").Value = 10.894872
$ws.Range("AN195").Value = 21.1051047
$ws.Range("AO195").Value = 28.0096818
$ws.Range("AP195").Value = 15.039184
$ws.Range("AQ195").Value = 14.2700365
$ws.Range("AS195").Value = 12.9001579
$ws.Range("AT195").Value = 27.249707
$ws.Range("AU195").Value = 23.8585882
$ws.Range("AV195").Value = 29.0621821
$ws.Range("AW195").Value = 27.8900001
$ws.Range("AX195").Value = 26.0379488
$ws.Range("AY195").Value = 18.2607489
$ws.Range("BA195").Value = 9.922699400000001
$ws.Range("BB195").Value = 17.5036242
$ws.Range("BC195").Value = 20.7147663
$ws.Range("BD195").Value = 21.1693121
$ws.Range("BE195").Value = 18.5597218

# --- Row 196 ---
$ws.Range("A196").Value = "13 08 2020"
$ws.Range("B196").Value = 18.6936937
$ws.Range("C196").Value = 33.7821909
$ws.Range("D196").Value = 28.0666772
$ws.Range("F196").Value = 23.3087473
$ws.Range("G196").Value = 18.6588514
$ws.Range("H196").Value = 16.8359336
$ws.Range("I196").Value = 9.821770799999999
$ws.Range("J196").Value = 12.71777
$ws.Range("K196").Value = 13.4557495
$ws.Range("L196").Value = 24.9574225
$ws.Range("M196").Value = 30.4019733
$ws.Range("O196").Value = 15.5840286
$ws.Range("P196").Value = 26.828655
$ws.Range("Q196").Value = 29.1936302
$ws.Range("R196").Value = 19.3952905
$ws.Range("S196").Value = 24.5529044
$ws.Range("T196").Value = 25.0408435
$ws.Range("U196").Value = 25.7780475
$ws.Range("V196").Value = 32.1366465
$ws.Range("W196").Value = 11.0280196
$ws.Range("X196").Value = 14.5178753
$ws.Range("Y196").Value = 11.0250298
$ws.Range("Z196").Value = 16.8503228
$ws.Range("AA196").Value = 19.8465747
$ws.Range("AB196").Value = 25.4309807
$ws.Range("AD196").Value = 34.6822104
$ws.Range("AE196").Value = 24.7782138
$ws.Range("AF196").Value = 20.2788918
$ws.Range("AG196").Value = 25.7583316
$ws.Range("AH196").Value = 25.6112411
$ws.Range("AI196").Value = 10.4150198
$ws.Range("AJ196").Value = 11.5306769
$ws.Range("AK196").Value = 17.5847567
$ws.Range("AL196").Value = 24.1996708
$ws.Range("AM196").Value = 10.9004385
$ws.Range("AN196").Value = 21.1175646
$ws.Range("AO196").Value = 28.2063465
$ws.Range("AP196").Value = 14.9578234
$ws.Range("AQ196").Value = 14.1791869
$ws.Range("AS196").Value = 12.3592479
$ws.Range("AT196").Value = 26.9412275
$ws.Range("AU196").Value = 24.1123584
$ws.Range("AV196").Value = 29.572785
$ws.Range("AW196").Value = 27.368372
$ws.Range("AX196").Value = 25.41985
$ws.Range("AY196").Value = 17.9813365
$ws.Range("BA196").Value = 10.2067241
$ws.Range("BB196").Value = 17.3355876
$ws.Range("BC196").Value = 20.751151
$ws.Range("BD196").Value = 20.3526514
$ws.Range("BE196").Value = 16.7832372

# --- Row 197 ---
$ws.Range("A197").Value = "14 08 2020"
$ws.Range("B197").Value = 19.1838649
$ws.Range("C197").Value = 33.427832
$ws.Range("D197").Value = 27.7531198
$ws.Range("F197").Value = 22.9382747
$ws.Range("G197").Value = 18.4053191
$ws.Range("H197").Value = 16.4702784
$ws.Range("I197").Value = 9.966548599999999
$ws.Range("J197").Value = 11.9930475
$ws.Range("K197").Value = 13.2533825
$ws.Range("L197").Value = 24.9258695
$ws.Range("M197").Value = 30.527199
$ws.Range("O197").Value = 15.3644315
$ws.Range("P197").Value = 25.560564
$ws.Range("Q197").Value = 29.4251373
$ws.Range("R197").Value = 19.3837299
$ws.Range("S197").Value = 24.4544311
$ws.Range("T197").Value = 25.3868491
$ws.Range("U197").Value = 26.174588
$ws.Range("V197").Value = 32.1311277
$ws.Range("W197").Value = 11.2599488
$ws.Range("X197").Value = 14.1990365
$ws.Range("Y197").Value = 10.617236
$ws.Range("Z197").Value = 16.7923574
$ws.Range("AA197").Value = 19.3419711
$ws.Range("AB197").Value = 25.6342516
$ws.Range("AD197").Value = 35.689749
$ws.Range("AE197").Value = 23.7015114
$ws.Range("AF197").Value = 19.9852792
$ws.Range("AG197").Value = 25.8311965
$ws.Range("AH197").Value = 25.8096258
$ws.Range("AI197").Value = 10.3065752
$ws.Range("AJ197").Value = 11.1431056
$ws.Range("AK197").Value = 16.9371142
$ws.Range("AL197").Value = 24.0018194
$ws.Range("AM197").Value = 10.7618898
$ws.Range("AN197").Value = 20.6417241
$ws.Range("AO197").Value = 27.2947301
$ws.Range("AP197").Value = 14.7829464
$ws.Range("AQ197").Value = 13.9146124
$ws.Range("AS197").Value = 12.0463044
$ws.Range("AT197").Value = 25.7603004
$ws.Range("AU197").Value = 24.6221963
$ws.Range("AV197").Value = 29.545132
$ws.Range("AW197").Value = 26.8449307
$ws.Range("AX197").Value = 25.5196967
$ws.Range("AY197").Value = 17.7430714
$ws.Range("BA197").Value = 9.2129619
$ws.Range("BB197").Value = 17.2592871
$ws.Range("BC197").Value = 20.3054107
$ws.Range("BD197").Value = 21.2220279
$ws.Range("BE197").Value = 18.9446464

# --- Row 198 ---
$ws.Range("A198").Value = "15 08 2020"
$ws.Range("B198").Value = 19.1987513
$ws.Range("C198").Value = 33.2747161
$ws.Range("D198").Value = 27.584159
$ws.Range("F198").Value = 22.0904279
$ws.Range("G198").Value = 18.3198168
$ws.Range("H198").Value = 16.6290883
$ws.Range("I198").Value = 9.743336599999999
$ws.Range("J198").Value = 11.7007673
$ws.Range("K198").Value = 12.7610966
$ws.Range("L198").Value = 24.907447
$ws.Range("M198").Value = 30.2053708
$ws.Range("O198").Value = 16.2385321
$ws.Range("P198").Value = 25.9850905
$ws.Range("Q198").Value = 29.6955876
$ws.Range("R198").Value = 19.3076758
$ws.Range("S198").Value = 24.4792813
$ws.Range("T198").Value = 25.9166873
$ws.Range("U198").Value = 26.1704314
$ws.Range("V198").Value = 32.1175371
$ws.Range("W198").Value = 11.1762899
$ws.Range("X198").Value = 14.3588571
$ws.Range("Y198").Value = 10.1107465
$ws.Range("Z198").Value = 16.967869
$ws.Range("AA198").Value = 19.2988243
$ws.Range("AB198").Value = 25.2583839
$ws.Range("AD198").Value = 34.9141754
$ws.Range("AE198").Value = 22.8502989
$ws.Range("AF198").Value = 19.9396592
$ws.Range("AG198").Value = 25.7258507
$ws.Range("AH198").Value = 25.5199562
$ws.Range("AI198").Value = 10.0619835
$ws.Range("AJ198").Value = 10.8683679
$ws.Range("AK198").Value = 17.2474982
$ws.Range("AL198").Value = 24.4235001
$ws.Range("AM198").Value = 10.5770344
$ws.Range("AN198").Value = 20.3573403
$ws.Range("AO198").Value = 27.1270502
$ws.Range("AP198").Value = 14.5561165
$ws.Range("AQ198").Value = 13.825994
$ws.Range("AS198").Value = 12.1341478
$ws.Range("AT198").Value = 25.6155314
$ws.Range("AU198").Value = 24.5177924
$ws.Range("AV198").Value = 28.9365391
$ws.Range("AW198").Value = 26.4887662
$ws.Range("AX198").Value = 25.3008116
$ws.Range("AY198").Value = 17.6560204
$ws.Range("BA198").Value = 9.374993399999999
$ws.Range("BB198").Value = 17.3321307
$ws.Range("BC198").Value = 20.2057439
$ws.Range("BD198").Value = 20.9229353
$ws.Range("BE198").Value = 18.2362417

# --- Row 199 ---
$ws.Range("A199").Value = "16 08 2020"
$ws.Range("B199").Value = 21.0907336
$ws.Range("C199").Value = 33.179098
$ws.Range("D199").Value = 28.5774491
$ws.Range("F199").Value = 21.8386044
$ws.Range("G199").Value = 17.9983606
$ws.Range("H199").Value = 16.4356827
$ws.Range("I199").Value = 9.318960300000001
$ws.Range("J199").Value = 11.8757613
$ws.Range("K199").Value = 12.3959001
$ws.Range("L199").Value = 24.5730907
$ws.Range("M199").Value = 29.7404412
$ws.Range("O199").Value = 16.2545235
$ws.Range("P199").Value = 26.157528
$ws.Range("Q199").Value = 29.436943
$ws.Range("R199").Value = 19.1925035
$ws.Range("S199").Value = 24.3479693
$ws.Range("T199").Value = 25.531528
$ws.Range("U199").Value = 25.4914788
$ws.Range("V199").Value = 31.527509
$ws.Range("W199").Value = 11.1632779
$ws.Range("X199").Value = 14.5476222
$ws.Range("Y199").Value = 10.0712768
$ws.Range("Z199").Value = 16.5792526
$ws.Range("AA199").Value = 19.4229163
$ws.Range("AB199").Value = 25.5561234
$ws.Range("AD199").Value = 34.4801424
$ws.Range("AE199").Value = 23.5342366
$ws.Range("AF199").Value = 19.6570704
$ws.Range("AG199").Value = 27.6955728
$ws.Range("AH199").Value = 24.686475
$ws.Range("AI199").Value = 10.259247
$ws.Range("AJ199").Value = 10.7386476
$ws.Range("AK199").Value = 16.8077537
$ws.Range("AL199").Value = 24.2840813
$ws.Range("AM199").Value = 10.4717989
$ws.Range("AN199").Value = 20.2996848
$ws.Range("AO199").Value = 26.7780441
$ws.Range("AP199").Value = 14.7600907
$ws.Range("AQ199").Value = 13.928179
$ws.Range("AS199").Value = 12.3523279
$ws.Range("AT199").Value = 25.0034232
$ws.Range("AU199").Value = 24.4485212
$ws.Range("AV199").Value = 29.2547936
$ws.Range("AW199").Value = 25.9248356
$ws.Range("AX199").Value = 23.9427531
$ws.Range("AY199").Value = 17.7829286
$ws.Range("BA199").Value = 8.717712199999999
$ws.Range("BB199").Value = 17.2681214
$ws.Range("BC199").Value = 20.4625111
$ws.Range("BD199").Value = 21.2966108
$ws.Range("BE199").Value = 18.5370784

# --- Row 200 ---
$ws.Range("A200").Value = "17 08 2020"
